{"js": "// Commit: \"Modified A2 design note\"\n//\n// 1) \"Interface Document along with the method ...\" -> remove the word\n//    \"along \" so it reads \"Interface Document with the method ...\".\n// 2) Remove the whole bullet item that discusses TreeSet being accessed\n//    by multiple threads (the \"If the TreeSet is accessed by multiple\n//    threads concurrently, ... it will be synchronized externally.\"\n//    paragraph).\n\nconst body = context.document.body;\n\n// --- 1) \"along with\" -> \"with\" --------------------------------------\nconst alongResults = body.search(\"Interface Document along with the method\", {\n  matchCase: true\n});\nalongResults.load(\"items\");\nawait context.sync();\n\nif (alongResults.items.length > 0) {\n  alongResults.items[0].insertText(\n    \"Interface Document with the method\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- 2) remove the \"accessed by multiple threads\" bullet paragraph --\nconst threadResults = body.search(\n  \"If the TreeSet is accessed by multiple threads concurrently\",\n  { matchCase: true }\n);\nthreadResults.load(\"items\");\nawait context.sync();\n\nif (threadResults.items.length > 0) {\n  const para = threadResults.items[0].paragraphs.getFirst();\n  para.delete();\n  await context.sync();\n}\n", "ps1": "# Commit: \"Modified A2 design note\"\n#\n# 1) \"Interface Document along with the method ...\" -> remove the word\n#    \"along \" so it reads \"Interface Document with the method ...\".\n# 2) Remove the whole bullet item that discusses TreeSet being accessed\n#    by multiple threads (the \"If the TreeSet is accessed by multiple\n#    threads concurrently, ... it will be synchronized externally.\"\n#    paragraph).\n\n$d = $word.ActiveDocument\n\n# --- 1) \"along with\" -> \"with\" ---------------------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"Interface Document along with the method\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"Interface Document with the method\"\n$find1.Forward = $true\n$find1.Wrap = 1          # wdFindContinue\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)  # wdReplaceAll\n\n# --- 2) remove the \"accessed by multiple threads\" bullet paragraph ---\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = \"If the TreeSet is accessed by multiple threads concurrently\"\n$find2.Forward = $true\n$find2.Wrap = 1          # wdFindContinue\n$found2 = $find2.Execute()\n\nif ($found2) {\n    $range2.Expand(4)   # wdParagraph\n    $range2.Delete()\n}\n"}
